$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename / restructure the header row ---
# Old headers: No | Kode Barang | Nama Barang | Harga Beli | Harga Jual | Kategori
# New headers: Kategori_ID | Barang_Kode | Barang_Nama | Harga_Beli | Harga_Jual
$ws.Range("A1").Value = "Kategori_ID"
$ws.Range("B1").Value = "Barang_Kode"
$ws.Range("C1").Value = "Barang_Nama"
$ws.Range("D1").Value = "Harga_Beli"
$ws.Range("E1").Value = "Harga_Jual"

# Drop the old 6th column ("Kategori") - the sheet is now only A:E
$ws.Range("F1").Clear()

# --- Add a second (blank/formatted) data row under the new B/C headers ---
$ws.Range("B2:C2").Font.Family = 0
$ws.Range("B2:C2").Font.Name = "Calibri"
$ws.Range("B2:C2").Font.Size = 11
$ws.Range("B2:C2").Font.Bold = $false
$ws.Range("B2:C2").Font.Color = 0x000000

# --- Column A is now wider (was the narrow "No" column) ---
$ws.Columns.Item(1).ColumnWidth = 15.57

# --- Selection / active cell moved away from the old data range ---
$ws.Range("K10").Select()
